$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.530.26'
$ws.Range('D3').Value = '1.842.92'
$ws.Range('E3').Value = '  +3.93%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.21'
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.563'
$ws.Range('E6').Value = '  +2.81%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.50'
$ws.Range('E8').Value = '  +3.75%  '
$ws.Range('E9').Value = '  +6.13%  '
$ws.Range('E10').Value = '  +10.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0934'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').Value = '2.108.76'
$ws.Range('E12').Value = '  +3.95%  '
$ws.Range('D13').Value = '1.837.74'
$ws.Range('E13').Value = '  +3.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.01'
$ws.Range('E14').Value = '  +1.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.652'
$ws.Range('E15').Value = '  +5.08%  '
$ws.Range('D16').Value = '34.582.20'
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.40'
$ws.Range('E17').Value = '  +5.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.21'
$ws.Range('E18').Value = '  +2.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '253.66'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('E20').Value = '  +10.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.31'
$ws.Range('E21').Value = '  +9.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('E23').Value = '  +3.28%  '
$ws.Range('E24').Value = '  +1.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.01'
$ws.Range('E25').Value = '  +3.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.87'
$ws.Range('E26').Value = '  +3.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.31'
$ws.Range('E27').Value = '  +4.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.116'
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').Value = '  +5.96%  '
$ws.Range('E31').Value = '  +2.36%  '
$ws.Range('B32').Value = 'Swop.fi'
$ws.Range('C32').Value = 'https://coinranking.com/coin/yrCr2HW2c+swopfi-swop'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '526.34'
$ws.Range('E32').Value = '  +910.33%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.22'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.65'
$ws.Range('E34').Value = '  +2.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.97'
$ws.Range('E35').Value = '  +7.39%  '
$ws.Range('D36').Value = '1.470.46'
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.659'
$ws.Range('E37').Value = '  +5.72%  '
$ws.Range('E38').Value = '  +3.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.984'
$ws.Range('E40').Value = '  +11.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '83.21'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('E44').Value = '  +6.37%  '
$ws.Range('E45').Value = '  +7.30%  '
$ws.Range('D46').Value = '2.006.01'
$ws.Range('E46').Value = '  +4.10%  '
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0496'
$ws.Range('E48').Value = '  -2.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '12.39'
$ws.Range('E49').Value = '  +5.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '106.83'
$ws.Range('E50').Value = '  +9.81%  '
$ws.Range('E51').Value = '  +0.15%  '
